# Update calibration estimates (moving-average std revisions) across the
# riskDriverType/AR/SETAR/GARCH/TARCH/AR_TARCH sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.00499731415748183
$ws.Range("B3").Value = 0.7711701032894739
$ws.Range("B4").Value = 0.1035095332260964
$ws.Range("B5").Value = "[0.9999999999999999, 0.025119562714767167, -0.04043037783277752, -0.11377522093295238, -0.1263982703342901, 0.05729696109786486, 0.011154612049982467, 0.07035491216299222, 0.03186409641904842, -0.012453649374926748, -0.06997161244419883, -0.04382851311469033, 0.0010308417392698084, 0.02886426710454328, 0.07250179301358817, 0.03522460735072164, -0.007025144569549078, -0.05583638288337831, -0.05284521720960309, 0.0011761930144414331]"

# ---------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1854976339519105
$ws.Range("B5").Value = 0.5547653944959511
$ws.Range("B6").Value = 0.06213472818997952
$ws.Range("B7").Value = 0.1724218551567645
$ws.Range("B8").Value = 0.5685865772672786
$ws.Range("B9").Value = 0.05971121150875135
$ws.Range("B10").Value = "[1.0, 0.09029284407118399, 0.003090160991593748, -0.01458868333275319, -0.013131832749729094, -0.013045224533461158, 0.00013197782382329455, 0.0018329649223643868, 0.0029397216042974334, 0.008145383099990983, -0.006604173561337797, -0.0014437052135839385, -0.0031953689725195025, 0.025676938455368777, 0.028010908412673285, 0.014193155799411782, 0.000842049264462496, -0.004714149521669019, 0.009819455021543506, -0.015215927179139159]"

# ---------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = [double]"7.078375205828765E-05"
$ws.Range("B3").Value = [double]"0.0001887705531063615"
$ws.Range("B4").Value = 0.003958963368569435
$ws.Range("B5").Value = 0.9943866867833515
$ws.Range("B6").Value = "[1.0, 0.027629632874653715, -0.030925524881919863, -0.1275375991688263, -0.14662833613375903, 0.08505365402715732, 0.0034969294511322984, 0.08532621533588629, 0.036943154463400636, -0.022910206460886057, -0.07214021264858068, -0.051896719295105045, 0.0016369591632326031, 0.036226090508683896, 0.08224966455505008, 0.02504311210321133, -0.014994620056548232, -0.06652901305197097, -0.06988438051247067, 0.0004693468076000554]"

# ---------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.0007433445166432991
$ws.Range("B3").Value = 0.1113250269049656
$ws.Range("B4").Value = 0.03984717846706533
$ws.Range("B5").Value = [double]"2.823159889204217E-13"
$ws.Range("B6").Value = "[1.0, -0.009140478439983042, -0.02638775392404581, -0.11844625024583522, -0.1450048596864767, 0.0939294527168299, 0.0007800027806599276, 0.08472041308542935, 0.037082593102983896, -0.02047629876139071, -0.06819600789608209, -0.047934913228888036, 0.0006323112185287519, 0.039943183840940304, 0.08580120702501767, 0.026619967649498807, -0.00948877996003385, -0.06391050196589737, -0.06542553312455232, 0.0072674302848237315]"
$ws.Range("B7").Value = 0.01199870732960074

# ---------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.005587142484813909
$ws.Range("B3").Value = 0.09947967699301188
$ws.Range("B4").Value = 0.04214799570264793
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[1.0, -0.012426916358293173, -0.039465188275681164, -0.11129039118049046, -0.12671586497594622, 0.06322222380553758, 0.007403567800449045, 0.06990618871533896, 0.03119693010874365, -0.010738379210647641, -0.06749166412294484, -0.04309687849464792, 0.0008267335248117876, 0.026313439589745344, 0.06988128392222186, 0.03405145663264218, -0.005893336838458587, -0.056225570131385544, -0.05300252652696254, 0.002034935012816862]"
$ws.Range("B7").Value = -0.008482083457821876
$ws.Range("B9").Value = 0.7704589490061118
